# Apply the "tabla de precios" update:
#  - Update labor cost (row 13) and network config cost (row 14)
#  - Insert a new "Capacitación" line item (new row 15), pushing the
#    Total row down to row 16 and updating its SUM range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update existing amounts
$ws.Range("J13").Value = 4000
$ws.Range("K13").Value = 4000

$ws.Range("J14").Value = 750
$ws.Range("K14").Value = 750

# 2) Insert a new row before the current Total row (row 15), shifting the
#    Total row down to row 16.
$ws.Rows.Item(15).Insert()

# Give the freshly-inserted row the same formatting as the Total row that
# just got pushed down to row 16 (matches the original row's look/borders
# instead of inheriting brand-new style entries from Insert()).
$ws.Range("D16:K16").Copy()
$ws.Range("D15:K15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3) Fill in the new "Capacitación" line (row 15)
$ws.Range("D15").Value = 9
$ws.Range("E15").Value = "Capacitación"
$ws.Range("H15").Value = "Proyecto"
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 600
$ws.Range("K15").Value = 600

# 4) Fix up the Total row formula (now row 16) to include the new row
$ws.Range("K16").Formula = "=SUM(K8:K15)"

# 5) Update the view state to match the saved workbook
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("G19").Select()
